$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 53

$ws.Cells.Item($row, 1).Value = "Globo"
$ws.Cells.Item($row, 2).Value = "RJ TV 2"
$ws.Cells.Item($row, 3).Value = "Obras"
$ws.Cells.Item($row, 4).Value = "2025-04-04T19:33"
$ws.Cells.Item($row, 5).Value = "Negativo"
$ws.Cells.Item($row, 6).Value = "Chuva na rodoviária de Campos. Pancadas de chuva trazem novos transtornos ao Shopping Estrada. Repórter *ao vivo* do local. Problema de infraestrutura do Shopping Estrada é recorrente. Pessoas não puderam sentar porque estava chovendo na parte coberta da rodoviária. Emissora recebeu alguns vídeos durante a chuva. Repórter questionou resposta da prefeitura para o problema, mas não teve resposta. "
